$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.000.60"
$ws.Range("E2").Value = "  -0.18%  "

$ws.Range("D3").Value = "1.630.01"
$ws.Range("E3").Value = "  -0.95%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").Value = "'214.27"
$ws.Range("E5").Value = "  -0.70%  "

$ws.Range("E6").Value = "  -0.80%  "

$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("E8").Value = "  -1.92%  "

$ws.Range("E9").Value = "  -3.19%  "

$ws.Range("E10").Value = "  -5.27%  "

$ws.Range("D11").Value = "'0.0788"
$ws.Range("E11").Value = "  -1.24%  "

$ws.Range("D12").Value = "1.857.21"
$ws.Range("E12").Value = "  -0.95%  "

$ws.Range("D13").Value = "1.625.55"
$ws.Range("E13").Value = "  +0.70%  "

$ws.Range("E14").Value = "  -1.96%  "

$ws.Range("E15").Value = "  -3.06%  "

$ws.Range("D16").Value = "26.013.03"
$ws.Range("E16").Value = "  -0.20%  "

$ws.Range("D17").Value = "0.0₃0741"
$ws.Range("E17").Value = "  -2.92%  "

$ws.Range("D18").Value = "'61.57"
$ws.Range("E18").Value = "  -3.10%  "

$ws.Range("E19").Value = "  -0.16%  "

$ws.Range("D20").Value = "'193.29"
$ws.Range("E20").Value = "  -0.63%  "

$ws.Range("D21").Value = "'4.26"
$ws.Range("E21").Value = "  -2.41%  "

$ws.Range("E22").Value = "  -3.84%  "

$ws.Range("D23").Value = "'6.07"
$ws.Range("E23").Value = "  -2.23%  "

$ws.Range("E24").Value = "  +1.70%  "

$ws.Range("D25").Value = "'144.30"
$ws.Range("E25").Value = "  +0.21%  "

$ws.Range("E26").Value = "  -0.10%  "

$ws.Range("E27").Value = "  -3.88%  "

$ws.Range("E28").Value = "  -2.35%  "

$ws.Range("D29").Value = "'15.32"
$ws.Range("E29").Value = "  -1.26%  "

$ws.Range("D30").Value = "'1.24"
$ws.Range("E30").Value = "  -0.91%  "

$ws.Range("E31").Value = "  -2.57%  "

$ws.Range("E32").Value = "  -4.30%  "

$ws.Range("E33").Value = "  -4.83%  "

$ws.Range("E34").Value = "  -3.22%  "

$ws.Range("E35").Value = "  -2.32%  "

$ws.Range("D36").Value = "1.119.92"
$ws.Range("E36").Value = "  -1.07%  "

$ws.Range("D37").Value = "'0.853"
$ws.Range("E37").Value = "  -5.84%  "

$ws.Range("E38").Value = "  -1.25%  "

$ws.Range("E39").Value = "  -3.60%  "

$ws.Range("D40").Value = "'0.0153"
$ws.Range("E40").Value = "  -2.47%  "

$ws.Range("D41").Value = "'98.33"
$ws.Range("E41").Value = "  -0.73%  "

$ws.Range("D42").Value = "'0.770"
$ws.Range("E42").Value = "  -3.46%  "

$ws.Range("D43").Value = "1.766.67"
$ws.Range("E43").Value = "  -0.97%  "

$ws.Range("E44").Value = "  -5.68%  "

$ws.Range("E45").Value = "  -1.82%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'54.59"
$ws.Range("E46").Value = "  -3.54%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.0529"
$ws.Range("E47").Value = "  +1.28%  "

$ws.Range("E48").Value = "  -0.53%  "

$ws.Range("D49").Value = "'0.414"
$ws.Range("E49").Value = "  -0.18%  "

$ws.Range("D50").Value = "'7.51"
$ws.Range("E50").Value = "  -3.68%  "

$ws.Range("D51").Value = "'1.01"
$ws.Range("E51").Value = "  +0.11%  "
